# Update the account-statement (Estado de Cuenta) table on Hoja1.
# The rows B16:J28 contain a list of (DocType, DocNumber, Name, Period,
# ValorMora, SalarioBasico) records. The underlying data was refreshed:
# previous periods were removed and new ones added, and the rows are now
# grouped by worker with periods in descending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Each entry: Row, DocNumber, Name, Period, ValorMora, SalarioBasico
$data = @(
    @(16, "1082856017", "RUBEN DARIO CUCUNUBA SALINAS",     "2303", 46400, 1160000),
    @(17, "1082856017", "RUBEN DARIO CUCUNUBA SALINAS",     "2302", 46400, 1160000),
    @(18, "1082856017", "RUBEN DARIO CUCUNUBA SALINAS",     "2301", 41760, 1160000),
    @(19, "1050944608", "ALBEIRO ENRIQUE PEREZ CAUSIL",     "2303", 34027, 1160000),
    @(20, "1007738481", "LUIS ANTONIO DE AVILA HERNANDEZ",  "2301", 40000, 1000000),
    @(21, "1007738481", "LUIS ANTONIO DE AVILA HERNANDEZ",  "2212",  6667, 1000000),
    @(22, "1050967361", "JUAN CAMILO HERNANDEZ VIGGIANI",   "2305", 37120, 1000000),
    @(23, "1050967361", "JUAN CAMILO HERNANDEZ VIGGIANI",   "2304", 46400, 1000000),
    @(24, "1047376325", "JUAN DAVID PAEZ CORTECERO",        "2303", 46400, 1160000),
    @(25, "1047376325", "JUAN DAVID PAEZ CORTECERO",        "2302", 46400, 1160000),
    @(26, "1047376325", "JUAN DAVID PAEZ CORTECERO",        "2301", 27840, 1160000),
    @(27, "1066734978", "JHON DEIVY GARCIA SANCHEZ",        "2305", 37120, 1160000),
    @(28, "1066734978", "JHON DEIVY GARCIA SANCHEZ",        "2304", 46400, 1160000)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $docNumber = $entry[1]
    $name = $entry[2]
    $period = $entry[3]
    $valorMora = $entry[4]
    $salarioBasico = $entry[5]

    $ws.Cells.Item($row, 3).Value = $docNumber
    $ws.Cells.Item($row, 4).Value = $name
    $ws.Cells.Item($row, 5).Value = $period
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
